$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '68.793.98'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.26%  '

$ws.Range("D3").Value = "'" + '3.312.71'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.39%  '

$ws.Range("D5").Value = "'" + '590.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.49%  '

$ws.Range("D6").Value = "'" + '186.48'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.56%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = "'" + '0.605'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.10%  '

$ws.Range("E9").Value = '  +5.24%  '

$ws.Range("D10").Value = "'" + '6.73'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.32%  '

$ws.Range("D11").Value = "'" + '0.423'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.98%  '

$ws.Range("D12").Value = "'" + '3.877.23'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.26%  '

$ws.Range("E13").Value = '  +0.39%  '

$ws.Range("D14").Value = "'" + '29.08'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.97%  '

$ws.Range("D15").Value = "'" + '68.817.47'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.47%  '

$ws.Range("D16").Value = "'" + '0.0000174'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.09%  '

$ws.Range("D17").Value = "'" + '3.297.23'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.88%  '

$ws.Range("D18").Value = "'" + '5.92'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.28%  '

$ws.Range("D19").Value = "'" + '13.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.56%  '

$ws.Range("D20").Value = "'" + '386.38'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.51%  '

$ws.Range("D21").Value = "'" + '7.84'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.73%  '

$ws.Range("D22").Value = "'" + '71.78'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.00%  '

$ws.Range("D23").Value = "'" + '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.12%  '

$ws.Range("E24").Value = '  +4.27%  '

$ws.Range("D25").Value = "'" + '0.520'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.65%  '

$ws.Range("E26").Value = '  +7.93%  '

$ws.Range("D27").Value = "'" + '9.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.07%  '

$ws.Range("D28").Value = "'" + '0.998'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.59%  '

$ws.Range("D29").Value = "'" + '5.92'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.15%  '

$ws.Range("D30").Value = "'" + '2.02'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.47%  '

$ws.Range("D31").Value = "'" + '1.34'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.44%  '

$ws.Range("D32").Value = "'" + '23.18'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.87%  '

$ws.Range("D33").Value = "'" + '7.27'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.21%  '

$ws.Range("E34").Value = '  +0.02%  '

$ws.Range("D35").Value = "'" + '1.56'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.93%  '

$ws.Range("D36").Value = "'" + '163.36'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.31%  '

$ws.Range("D37").Value = "'" + '1.90'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.51%  '

$ws.Range("D38").Value = "'" + '0.842'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.99%  '

$ws.Range("D39").Value = "'" + '27.04'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.80%  '

$ws.Range("E40").Value = '  -0.69%  '

$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = "'" + '4.65'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.17%  '

$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").Value = "'" + '2.67'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.86%  '

$ws.Range("D43").Value = "'" + '25.91'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.37%  '

$ws.Range("D44").Value = "'" + '0.0699'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.10%  '

$ws.Range("D45").Value = "'" + '41.42'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.66%  '

$ws.Range("D46").Value = "'" + '2.652.03'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.86%  '

$ws.Range("D47").Value = "'" + '343.79'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.15%  '

$ws.Range("D48").Value = "'" + '0.0287'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.51%  '

$ws.Range("D49").Value = "'" + '32.48'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.65%  '

$ws.Range("E50").Value = '  +2.12%  '

$ws.Range("E51").Value = '  +0.72%  '
